$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.004.28"
$ws.Range("E2").Value = "  +5.69%  "
$ws.Range("D3").Value = "1.911.69"
$ws.Range("E3").Value = "  +4.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9989"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.61%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4728"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.22%  "
$ws.Range("E8").Value = "  +6.71%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.024"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +7.25%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.095"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.68%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.888.48"
$ws.Range("E14").Value = "  +3.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.381"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.15%  "
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001052"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06633"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9985"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.59%  "
$ws.Range("D22").Value = "29.043.43"
$ws.Range("E22").Value = "  +5.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.566"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.265"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").Value = "2.124.99"
$ws.Range("E26").Value = "  +4.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.91"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.76%  "
$ws.Range("E29").Value = "  +6.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.541"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.20%  "
$ws.Range("E31").Value = "  +2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.019"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09581"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.647"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.419"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.431"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06207"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02295"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.681"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.87%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.202"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6051"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.19%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1905"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.58"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9975"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.12%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5645"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.989"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07287"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.154"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +18.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.12"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.03%  "
